# Updates the cryptos worksheet with latest scraped price/volume data.
# Applies the GitHub Actions commit: "Updated cryptos list ... with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (prices, % volumes, and two swapped coin rows).
$updates = @{
    'D2' = '27.144.51'
    'E2' = '  +0.04%  '
    'D3' = '1.832.81'
    'E3' = '  +0.21%  '
    'E4' = '  -0.05%  '
    'D5' = '312.37'
    'E5' = '  -0.11%  '
    'E6' = '  -0.02%  '
    'E7' = '  -1.51%  '
    'D8' = '0.3706'
    'E8' = '  +1.05%  '
    'D9' = '0.07362'
    'E9' = '  -0.46%  '
    'D10' = '0.8738'
    'E10' = '  -0.79%  '
    'D11' = '0.07970'
    'E11' = '  +3.18%  '
    'D12' = '19.90'
    'E12' = '  -1.97%  '
    'D13' = '1.805.46'
    'E13' = '  -5.39%  '
    'B14' = 'Chainlink'
    'C14' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D14' = '6.594'
    'E14' = '  +0.85%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '5.347'
    'E15' = '  -0.52%  '
    'D16' = '92.10'
    'E16' = '  -1.47%  '
    'E17' = '  +0.38%  '
    'D18' = '0.000008897'
    'E18' = '  +1.91%  '
    'D19' = '1.008'
    'E19' = '  +0.10%  '
    'B20' = 'WrappedBTC'
    'C20' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D20' = '27.418.26'
    'E20' = '  -0.84%  '
    'B21' = 'Avalanche'
    'C21' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D21' = '14.72'
    'E21' = '  +0.60%  '
    'D22' = '5.143'
    'E22' = '  -1.96%  '
    'D23' = '10.62'
    'E23' = '  -0.05%  '
    'D24' = '2.009.70'
    'E24' = '  -3.65%  '
    'D25' = '152.51'
    'E25' = '  +1.00%  '
    'D26' = '1.839'
    'E26' = '  -1.77%  '
    'D27' = '18.54'
    'E27' = '  +0.20%  '
    'D28' = '2.097'
    'E28' = '  -2.11%  '
    'D29' = '5.094'
    'E29' = '  -1.82%  '
    'D30' = '115.60'
    'E30' = '  -0.84%  '
    'E31' = '  -0.56%  '
    'D32' = '2.969'
    'E32' = '  +1.02%  '
    'D33' = '0.7359'
    'E33' = '  -1.33%  '
    'E34' = '  -1.32%  '
    'D35' = '1.140'
    'E35' = '  -2.22%  '
    'D36' = '2.469'
    'E36' = '  -5.32%  '
    'D37' = '1.077'
    'E37' = '  -1.27%  '
    'E38' = '  +0.60%  '
    'D39' = '0.05244'
    'E39' = '  -1.22%  '
    'D40' = '2.937'
    'E40' = '  +0.21%  '
    'D41' = '7.184'
    'E41' = '  -2.01%  '
    'D42' = '0.5202'
    'E42' = '  -1.42%  '
    'D43' = '0.8672'
    'E43' = '  -13.95%  '
    'E44' = '  -0.56%  '
    'D45' = '8.254'
    'E45' = '  -1.79%  '
    'E46' = '  -1.36%  '
    'D47' = '10.24'
    'E47' = '  -1.90%  '
    'D48' = '1.007'
    'E48' = '  +0.03%  '
    'D49' = '102.68'
    'E49' = '  -1.80%  '
    'E50' = '  -1.26%  '
    'D51' = '0.06237'
    'E51' = '  -0.66%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)

    # Force the cell to remain plain text so strings such as "19.90", "1.140" or
    # "0.07970" keep their exact character representation instead of being
    # reinterpreted (and reformatted) as numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}
